$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-18 Tuesday" "2025-02-19 Wednesday"

Replace-Text "50×20=" "38×58="
Replace-Text "97×77=" "27×56="
Replace-Text "31×21=" "70×99="
Replace-Text "92×69=" "84×47="
Replace-Text "52×56=" "45×75="

Replace-Text "53×41=" "28×17="
Replace-Text "79×85=" "89×48="
Replace-Text "24×92=" "77×30="
Replace-Text "48×32=" "77×68="
Replace-Text "31×53=" "27×65="

Replace-Text "55×34=" "52×32="
Replace-Text "20×31=" "89×35="
Replace-Text "13×16=" "96×23="
Replace-Text "24×30=" "27×86="
Replace-Text "34×23=" "11×63="

Replace-Text "27×82=" "62×68="
Replace-Text "42×92=" "66×47="
Replace-Text "80×64=" "54×96="
Replace-Text "62×60=" "80×11="
Replace-Text "47×97=" "50×50="

Replace-Text "14×86=" "60×62="
Replace-Text "98×98=" "98×52="
Replace-Text "32×66=" "88×85="
Replace-Text "51×43=" "81×89="
Replace-Text "90×39=" "18×75="
